$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 22: E22 used to be a formula (D22*B22); it is now replaced by a
#     part-number label "3083et" (plain text, no formula). ---
$ws.Range("E22").Value = "3083et"

# --- Row 25: D25 quantity/cost updated, E25 (=D25*B25) recalculates
#     automatically since it still holds its formula. ---
$ws.Range("D25").Value = 7.29

# --- New row 59: "Front Panel Keys" line item. ---
$ws.Range("A59").Value = "Front Panel Keys"
$ws.Range("B59").Value = 1
$ws.Range("C59").Value = 1
$ws.Range("D59").Value = 1
$ws.Range("E59").Value = 1
$ws.Range("F59").Value = "n/a"
$ws.Range("G59").Value = "n/a"
$ws.Range("H59").Value = "n/a"

# --- Update the view: scroll back to the top (removing the old
#     topLeftCell="A4" saved position) and move the selection to J2. ---
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("J2").Select()
